$wb = $excel.ActiveWorkbook

# Sheet "1" (2nd worksheet in the workbook) holds the EF-Core exam hints table
# that gets two new rows appended (Equals / GetHashCode overrides).
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- new column D (extra remark column used by the two new rows) ---
$ws2.Columns.Item(4).ColumnWidth = 37.42578125

# --- row 18 : override Equals ---
$ws2.Range("A18").Value = 17
$ws2.Range("B18").Value = "override Equals"
$ws2.Range("C18").Value = "За сравняване на обекти се правят двата OVERRIDE-a: Equals и GetHashCode"
$ws2.Range("D18").Value = @"
        public override bool Equals(object obj)
        {
            var other = (CategoryProduct)obj;
            return this.CategoryId == other.CategoryId &&
                this.ProductId == other.ProductId;
        }
"@

# --- row 19 : override GetHashCode ---
$ws2.Range("A19").Value = 18
$ws2.Range("B19").Value = "override GetHashCode"
$ws2.Range("C19").Value = "За сравняване на обекти се правят двата OVERRIDE-a: Equals и GetHashCode"
$ws2.Range("D19").Value = @"
        public override int GetHashCode()
        {
            int hash = 13;
            hash = (hash * 7) + CategoryId.GetHashCode();
            hash = (hash * 7) + ProductId.GetHashCode();
            return hash;
        }
"@

# --- formatting: reuse the look of the existing table rows/columns ---
# Column A (row numbers) -> same style as A15:A17
$ws2.Range("A15").Copy()
$ws2.Range("A18:A19").PasteSpecial(-4122)

# Columns B/C (title + shared explanation) -> same style as B15
$ws2.Range("B15").Copy()
$ws2.Range("B18:C19").PasteSpecial(-4122)

# Column D (code remark) -> same wrap-text style used for the other code cells (sheet "Лист1" A30)
$ws1.Range("A30").Copy()
$ws2.Range("D18:D19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws2.Rows.Item(18).RowHeight = 150
$ws2.Rows.Item(19).RowHeight = 165

# --- view state: scroll down and select B18 like the saved workbook ---
$ws2.Activate()
$ws2.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 16
